$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42 and 43 swap their Coin/Link/Price data, with distinct new Volume values
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.19%  "

$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.04%  "

# Remaining price (D) and volume (E) updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.849.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.317.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.85"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.313.00"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.577"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.74"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.849.00"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "569.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -9.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.830.61"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.321.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.64"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.888"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.99"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.74"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -8.49%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.33"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.40"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -9.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "556.28"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.81"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.16%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.744.52"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.52"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.60"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0₃0683"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.59%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -11.98%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.29%  "
